$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy formatting from column R (year 2022) into new column S (year 2023)
$ws.Range("R3:R14").Copy() | Out-Null
$ws.Range("S3:S14").PasteSpecial(-4122) | Out-Null

# Set the new values for column S
$ws.Range("S3").Value = 2023
$ws.Range("S4").Value = 18.4
$ws.Range("S5").Value = 18.4
$ws.Range("S6").Value = 493
$ws.Range("S7").Value = 465
$ws.Range("S8").Value = 797.6
$ws.Range("S9").Value = 10.8
$ws.Range("S10").Value = 4.5
$ws.Range("S11").Value = 7.6
$ws.Range("S12").Value = 12
$ws.Range("S13").Value = 9.2
$ws.Range("S14").Value = "_"

# Move the selection like in the target sheet
$ws.Range("D21").Select() | Out-Null
